$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# Row 16 already exists - update its Open points / Change log columns and
# clear the (now unused) pyinstaller version cell in column G.
# ---------------------------------------------------------------------------
$ws.Range("B16").Value = "AUTOMATA CELULAR - copia (25)"
$ws.Range("C16").Value = "-Make the code more readable using functions.`n-UI: Delete rows according to working functionality.`n*Graphic representation of F'.`n-Rework E calc.`n-Document every function.`n_OPTIONAL: Make it possible to reduce Niches on mutations.`n-When there is no data, the program crashes.`n-Check save data formatting.`n-Rework graphic representation to make an EXE."
$ws.Range("D16").Value = "-Fixed group selection.`n-Fixed reciprocal association.`n-Graphic representation disabled.`n-No more EXE from pyinstaller."
$ws.Range("G16").Value = ""

# ---------------------------------------------------------------------------
# Copy the formatting of the existing alternating rows (17 = odd style,
# 18 = even style) down to the rows we are about to populate / create
# (19 through 30), preserving the style banding already used in the sheet.
# ---------------------------------------------------------------------------
$ws.Range("A17:G17").Copy() | Out-Null
$ws.Range("A19:G19").PasteSpecial(-4122) | Out-Null
$ws.Range("A21:G21").PasteSpecial(-4122) | Out-Null
$ws.Range("A23:G23").PasteSpecial(-4122) | Out-Null
$ws.Range("A25:G25").PasteSpecial(-4122) | Out-Null
$ws.Range("A27:G27").PasteSpecial(-4122) | Out-Null
$ws.Range("A29:G29").PasteSpecial(-4122) | Out-Null

$ws.Range("A18:G18").Copy() | Out-Null
$ws.Range("A20:G20").PasteSpecial(-4122) | Out-Null
$ws.Range("A22:G22").PasteSpecial(-4122) | Out-Null
$ws.Range("A24:G24").PasteSpecial(-4122) | Out-Null
$ws.Range("A26:G26").PasteSpecial(-4122) | Out-Null
$ws.Range("A28:G28").PasteSpecial(-4122) | Out-Null
$ws.Range("A30:G30").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Row 17 - version 0.5.0
# ---------------------------------------------------------------------------
$ws.Range("A17").Value = "0.5.0"
$ws.Range("B17").Value = "AUTOMATA CELULAR - copia (27)"
$ws.Range("C17").Value = "-Make the code more readable using functions.`n-UI: Delete rows according to working functionality.`n-Document every function.`n_OPTIONAL: Make it possible to reduce Niches on mutations.`n-When there is no data, the program crashes.`n-Check save data formatting."
$ws.Range("D17").Value = "-Graphic representation enabled.`n-Graphic representation of F'.`n-Reworked E calc.`n-Reworked file system"
$ws.Range("E17").Value = "Python 3.6.1"
$ws.Range("F17").Value = "Qt version: 5.6.2`nSIP version: 4.18`nPyQt version: 5.6"
$ws.Range("G17").Value = ""

# ---------------------------------------------------------------------------
# Row 18 - version 0.5.2
# ---------------------------------------------------------------------------
$ws.Range("A18").Value = "0.5.2"
$ws.Range("B18").Value = "AUTOMATA CELULAR - copia (28)"
$ws.Range("C18").Value = "-Make the code more readable using functions.`n-UI: Delete rows according to working functionality.`n-Document every function.`n_OPTIONAL: Make it possible to reduce Niches on mutations.`n-When there is no data, the program crashes.`n-Check save data formatting."
$ws.Range("D18").Value = "-Reworked E calc.`n-Reworked P calc."
$ws.Range("E18").Value = "Python 3.6.1"
$ws.Range("F18").Value = "Qt version: 5.6.2`nSIP version: 4.18`nPyQt version: 5.6"
$ws.Range("G18").Value = ""

# ---------------------------------------------------------------------------
# Row 19 - version 0.6.0
# ---------------------------------------------------------------------------
$ws.Range("A19").Value = "0.6.0"
$ws.Range("B19").Value = "AUTOMATA CELULAR - copia (29)"
$ws.Range("C19").Value = "-Make the code more readable using functions.`n-UI: Delete rows according to working functionality.`n-Document every function.`n_OPTIONAL: Make it possible to reduce Niches on mutations.`n-When there is no data, the program crashes.`n-Check save data formatting."
$ws.Range("D19").Value = "-Reworked GS to be based on Greed`n-Results saved at the end for Gapminder"
$ws.Range("E19").Value = "Python 3.6.1"
$ws.Range("F19").Value = "Qt version: 5.6.2`nSIP version: 4.18`nPyQt version: 5.6"
$ws.Range("G19").Value = ""

# ---------------------------------------------------------------------------
# Row 20 - version 0.6.1
# ---------------------------------------------------------------------------
$ws.Range("A20").Value = "0.6.1"
$ws.Range("B20").Value = "AUTOMATA CELULAR - copia (30)"
$ws.Range("C20").Value = "-UI: Delete rows according to working functionality.`n-Document every function.`n_OPTIONAL: Make it possible to reduce Niches on mutations.`n-When there is no data, the program crashes."
$ws.Range("D20").Value = "-Reworked functions and unified them.`n-Finished no GUI option"
$ws.Range("E20").Value = "Python 3.6.1"
$ws.Range("F20").Value = "Qt version: 5.6.2`nSIP version: 4.18`nPyQt version: 5.6"
$ws.Range("G20").Value = ""

# ---------------------------------------------------------------------------
# Row 21 - version 0.6.2
# ---------------------------------------------------------------------------
$ws.Range("A21").Value = "0.6.2"
$ws.Range("B21").Value = "AUTOMATA CELULAR - copia (31)"
$ws.Range("C21").Value = "-UI: Delete rows according to working functionality.`n-Document every function.`n_OPTIONAL: Make it possible to reduce Niches on mutations.`n-When there is no data, the program crashes."
$ws.Range("D21").Value = "-Reworked grouping`n-Fixed missundertanding for Individual Selection Pressure.`n-GUI: Group Selection deaths percentage made float from int"
$ws.Range("E21").Value = "Python 3.6.1"
$ws.Range("F21").Value = "Qt version: 5.6.2`nSIP version: 4.18`nPyQt version: 5.6"
$ws.Range("G21").Value = ""

# ---------------------------------------------------------------------------
# Rows 22-30 stay empty (blank spares at the bottom of the table), only the
# formatting copied above is needed for them.
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# Row heights for the newly documented rows.
# ---------------------------------------------------------------------------
$ws.Rows.Item(17).RowHeight = 86.4
$ws.Rows.Item(18).RowHeight = 86.4
$ws.Rows.Item(19).RowHeight = 86.4
$ws.Rows.Item(20).RowHeight = 57.6
$ws.Rows.Item(21).RowHeight = 72

# ---------------------------------------------------------------------------
# Update the view: scroll so column D is leftmost, keep the header row
# frozen, and select D22 (the first empty spare row).
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("D1").Select()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("D22").Select()
